$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) values: force text type (NumberFormat '@') so numeric-looking
# strings like '83.69' are not coerced to floating-point numbers, then reset the
# cell style back to Normal so no stray style index is left behind.
function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "39.669.45"
$ws.Range("E2").Value = "  -1.01%  "
Set-TextValue $ws "D3" "2.216.31"
$ws.Range("E3").Value = "  -5.37%  "
$ws.Range("E4").Value = "  +0.01%  "
Set-TextValue $ws "D5" "298.30"
$ws.Range("E5").Value = "  -4.06%  "
Set-TextValue $ws "D6" "83.69"
$ws.Range("E6").Value = "  -2.07%  "
$ws.Range("E7").Value = "  -2.92%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -3.89%  "
Set-TextValue $ws "D10" "0.0783"
$ws.Range("E10").Value = "  -3.34%  "
Set-TextValue $ws "D11" "29.67"
$ws.Range("E11").Value = "  -1.29%  "
Set-TextValue $ws "D12" "46.10"
$ws.Range("E12").Value = "  -12.12%  "
Set-TextValue $ws "D14" "2.556.13"
$ws.Range("E14").Value = "  -5.30%  "
$ws.Range("E15").Value = "  -2.41%  "
Set-TextValue $ws "D16" "14.14"
$ws.Range("E16").Value = "  -4.07%  "
Set-TextValue $ws "D17" "2.221.14"
$ws.Range("E17").Value = "  -6.17%  "
Set-TextValue $ws "D18" "0.719"
$ws.Range("E18").Value = "  -5.39%  "
Set-TextValue $ws "D19" "39.587.96"
$ws.Range("E19").Value = "  -1.11%  "
Set-TextValue $ws "D20" "0.0₃0880"
$ws.Range("E20").Value = "  -2.76%  "
Set-TextValue $ws "D21" "5.75"
$ws.Range("E21").Value = "  -6.02%  "
Set-TextValue $ws "D22" "65.06"
$ws.Range("E22").Value = "  -4.58%  "
Set-TextValue $ws "D23" "10.43"
$ws.Range("E23").Value = "  -2.65%  "
Set-TextValue $ws "D24" "232.47"
$ws.Range("E24").Value = "  -1.23%  "
$ws.Range("E25").Value = "  -0.10%  "
Set-TextValue $ws "D26" "2.42"
$ws.Range("E26").Value = "  -5.15%  "
Set-TextValue $ws "D27" "1.83"
$ws.Range("E27").Value = "  +0.68%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws "D28" "2.30"
$ws.Range("E28").Value = "  +7.88%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws "D29" "22.75"
$ws.Range("E29").Value = "  -2.59%  "
Set-TextValue $ws "D30" "9.19"
$ws.Range("E30").Value = "  -1.20%  "
Set-TextValue $ws "D31" "32.27"
$ws.Range("E31").Value = "  -6.95%  "
$ws.Range("E32").Value = "  -2.61%  "
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("E34").Value = "  -5.48%  "
$ws.Range("E35").Value = "  -2.76%  "
Set-TextValue $ws "D36" "0.0703"
$ws.Range("E36").Value = "  -2.48%  "
Set-TextValue $ws "D37" "16.12"
$ws.Range("E37").Value = "  +3.00%  "
$ws.Range("E38").Value = "  -2.45%  "
$ws.Range("E39").Value = "  -1.58%  "
Set-TextValue $ws "D40" "2.66"
$ws.Range("E40").Value = "  -5.91%  "
$ws.Range("E41").Value = "  -4.63%  "
$ws.Range("E42").Value = "  -5.27%  "
Set-TextValue $ws "D43" "1.930.04"
$ws.Range("E43").Value = "  -1.34%  "
$ws.Range("E44").Value = "  -3.36%  "
$ws.Range("E45").Value = "  +0.66%  "
Set-TextValue $ws "D46" "9.23"
$ws.Range("E46").Value = "  -1.68%  "
Set-TextValue $ws "D47" "16.19"
$ws.Range("E47").Value = "  -8.39%  "
$ws.Range("E48").Value = "  -3.82%  "
Set-TextValue $ws "D49" "2.429.00"
$ws.Range("E49").Value = "  -5.10%  "
Set-TextValue $ws "D50" "71.08"
$ws.Range("E50").Value = "  +0.56%  "
Set-TextValue $ws "D51" "88.82"
$ws.Range("E51").Value = "  -4.51%  "
